# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K" - strikeouts) values for rows 2-26 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 7
    4  = 3
    5  = 3
    6  = 3
    7  = 5
    8  = 4
    9  = 2
    10 = 3
    11 = 4
    12 = 4
    13 = 2
    14 = 0
    15 = 5
    16 = 5
    17 = 3
    18 = 5
    19 = 7
    20 = 5
    21 = 4
    22 = 5
    23 = 5
    24 = 2
    25 = 5
    26 = 7
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
